$d = $word.ActiveDocument

# Original paragraph text is "Version 2." built from runs:
#   "Versi" | "on" | " 2" | "." (with a bookmark between " 2" and ".")
# Target paragraph text is "Version 1." built from runs:
#   "Version" | " 1." (bookmark kept, trailing "." run removed)

# Step 1: merge "Versi" + "on" into a single "Version" run.
# Restricting Find to the exact sub-range avoids touching the rest of the
# paragraph / merging unrelated runs.
$r1 = $d.Range(0, 7)
$r1.Find.Execute("Version", $true, $false, $false, $false, $false, $true, 1, $false, "Version", 2)

# Step 2: change the " 2" run into " 1."
$r2 = $d.Range(7, 9)
$r2.Find.Execute(" 2", $true, $false, $false, $false, $false, $true, 1, $false, " 1.", 2)

# Step 3: delete the now-trailing "." run entirely.
$full = $d.Content.Text
$r3 = $d.Range($full.Length - 2, $full.Length - 1)
$r3.Delete()
